# Update column F (dSF) values on Sheet1 to reflect repulled data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = 1
    "F5"  = 2
    "F7"  = -6
    "F8"  = -7
    "F9"  = -3
    "F10" = 0
    "F13" = -5
    "F14" = 17
    "F15" = -12
    "F16" = -2
    "F20" = -1
    "F21" = 4
    "F25" = 0
    "F29" = 2
    "F31" = 0
    "F34" = 1
    "F36" = -3
    "F37" = -7
    "F45" = -2
    "F50" = -3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
